$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $text
}

Set-CellText 1 1 "61÷9=6, 7"
Set-CellText 1 2 "84÷5=16, 4"
Set-CellText 1 3 "39÷6=6, 3"
Set-CellText 1 4 "49÷4=12, 1"
Set-CellText 1 5 "42÷5=8, 2"

Set-CellText 5 1 "54÷2=27, 0"
Set-CellText 5 2 "51÷7=7, 2"
Set-CellText 5 3 "97÷9=10, 7"
Set-CellText 5 4 "37÷7=5, 2"
Set-CellText 5 5 "26÷3=8, 2"

Set-CellText 9 1 "99÷6=16, 3"
Set-CellText 9 2 "97÷3=32, 1"
Set-CellText 9 3 "43÷5=8, 3"
Set-CellText 9 4 "94÷8=11, 6"
Set-CellText 9 5 "98÷8=12, 2"

Set-CellText 13 1 "13÷9=1, 4"
Set-CellText 13 2 "63÷3=21, 0"
Set-CellText 13 3 "53÷9=5, 8"
Set-CellText 13 4 "83÷5=16, 3"
Set-CellText 13 5 "97÷9=10, 7"

Set-CellText 17 1 "59÷4=14, 3"
Set-CellText 17 2 "14÷9=1, 5"
Set-CellText 17 3 "98÷2=49, 0"
Set-CellText 17 4 "54÷6=9, 0"
Set-CellText 17 5 "22÷5=4, 2"
